$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.240.23'
$ws.Range("E2").Value = '  +1.50%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.026.66'
$ws.Range("E3").Value = '  +3.35%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.75'
$ws.Range("E5").Value = '  +1.29%  '
$ws.Range("E6").Value = '  +1.97%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.22'
$ws.Range("E7").Value = '  -2.30%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +4.70%  '
$ws.Range("E10").Value = '  +2.13%  '
$ws.Range("E11").Value = '  +2.28%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.22'
$ws.Range("E12").Value = '  +6.71%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.61'
$ws.Range("E13").Value = '  +1.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.860'
$ws.Range("E14").Value = '  +3.13%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.320.00'
$ws.Range("E15").Value = '  +3.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.52'
$ws.Range("E16").Value = '  +4.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.029.57'
$ws.Range("E17").Value = '  +3.81%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.148.05'
$ws.Range("E18").Value = '  +1.58%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.70'
$ws.Range("E19").Value = '  +1.33%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0868'
$ws.Range("E20").Value = '  +1.58%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.27'
$ws.Range("E21").Value = '  +3.69%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '231.16'
$ws.Range("E22").Value = '  +0.40%  '
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("E24").Value = '  +2.76%  '
$ws.Range("E25").Value = '  +0.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.49'
$ws.Range("E26").Value = '  +3.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '163.87'
$ws.Range("E27").Value = '  +2.06%  '
$ws.Range("E28").Value = '  -3.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.84'
$ws.Range("E29").Value = '  +2.18%  '
$ws.Range("E30").Value = '  +7.53%  '
$ws.Range("E31").Value = '  +2.09%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.85'
$ws.Range("E32").Value = '  +1.63%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0666'
$ws.Range("E33").Value = '  +8.25%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.55'
$ws.Range("E34").Value = '  +2.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.48'
$ws.Range("E35").Value = '  +9.14%  '
$ws.Range("E36").Value = '  -3.33%  '
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.80'
$ws.Range("E38").Value = '  +1.68%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.44'
$ws.Range("E39").Value = '  -1.34%  '
$ws.Range("E40").Value = '  +0.26%  '
$ws.Range("E41").Value = '  +1.11%  '
$ws.Range("E42").Value = '  +1.46%  '
$ws.Range("E43").Value = '  +1.93%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.80'
$ws.Range("E44").Value = '  +4.60%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.33'
$ws.Range("E45").Value = '  +4.05%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.393.50'
$ws.Range("E46").Value = '  +1.71%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.07'
$ws.Range("E47").Value = '  +3.39%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.49'
$ws.Range("E48").Value = '  +4.94%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.20'
$ws.Range("E49").Value = '  +19.16%  '
$ws.Range("E50").Value = '  +0.38%  '
$ws.Range("E51").Value = '  +3.23%  '
